# Update sample names in column A so they match the seqtab naming
# convention (adds the dash-separated "-LBx" suffix and the Sx[/_L001]
# sample index used by the sequencer), required for script
# 08.phyloseq to perform correctly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value  = "052-LB1_S14"
$ws.Range("A3").Value  = "052-LB2_S15"
$ws.Range("A4").Value  = "052-LB3_S16"
$ws.Range("A5").Value  = "052-LB4_S17"
$ws.Range("A6").Value  = "052-LB5_S18"
$ws.Range("A7").Value  = "056-LB1_S19"
$ws.Range("A8").Value  = "056-LB2_S20"
$ws.Range("A9").Value  = "056-LB3_S21"
$ws.Range("A10").Value = "056-LB4_S22"
$ws.Range("A11").Value = "056-LB5_S23"
$ws.Range("A12").Value = "131-LB1_S24"
$ws.Range("A13").Value = "Xnc-LB1_S25"
$ws.Range("A14").Value = "MC_S26_L001"
$ws.Range("A15").Value = "PCRnc_S27_L"

# Matches the author's final cursor position recorded in the saved file.
$ws.Range("D14").Select()
